{"js": "// Apply the Gorilla Gold Megaways review text updates.\nconst replacements = [\n  [\n    \"Play Gorilla Gold Megaways Slot for Free - Review 2021\",\n    \"Play Gorilla Gold Megaways Free | Review & Gameplay\",\n  ],\n  [\n    \"Fun and exciting gameplay mechanics and features\",\n    \"Megaways engine with up to 470,000 ways to win\",\n  ],\n  [\n    \"Many pay lines and ways to win with the Megaways engine\",\n    \"Gorilla Gold Bonus with free spins and super spins\",\n  ],\n  [\n    \"Generous payouts and a high RTP value of 97%\",\n    \"Impressive RTP value of 97%\",\n  ],\n  [\n    \"Super Spins function requires a higher bet\",\n    \"Minimum bet limit of \\u20AC0.10 may not suit all players\",\n  ],\n  [\n    \"May be overwhelming for new players due to the multiple gaming zones\",\n    \"Limited number of bonus features\",\n  ],\n  [\n    \"Find out how to play Gorilla Gold Megaways online slot for free. Read our review to learn about the gameplay mechanics, features, payouts and more.\",\n    \"Discover the exciting gameplay of Gorilla Gold Megaways and play for free. Read our review for more details.\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [search, replacement] of replacements) {\n  const found = body.search(search, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(replacement, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-AllText($find, $replace) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n\nReplace-AllText \"Play Gorilla Gold Megaways Slot for Free - Review 2021\" \"Play Gorilla Gold Megaways Free | Review & Gameplay\"\nReplace-AllText \"Fun and exciting gameplay mechanics and features\" \"Megaways engine with up to 470,000 ways to win\"\nReplace-AllText \"Many pay lines and ways to win with the Megaways engine\" \"Gorilla Gold Bonus with free spins and super spins\"\nReplace-AllText \"Generous payouts and a high RTP value of 97%\" \"Impressive RTP value of 97%\"\nReplace-AllText \"Super Spins function requires a higher bet\" \"Minimum bet limit of \u20ac0.10 may not suit all players\"\nReplace-AllText \"May be overwhelming for new players due to the multiple gaming zones\" \"Limited number of bonus features\"\nReplace-AllText \"Find out how to play Gorilla Gold Megaways online slot for free. Read our review to learn about the gameplay mechanics, features, payouts and more.\" \"Discover the exciting gameplay of Gorilla Gold Megaways and play for free. Read our review for more details.\"\n"}
